# Update the two-digit / one-digit division problems in the worksheet table.
# Each old expression is unique within the document, so a straightforward
# Find/Replace (ReplaceAll) for each pair is safe. The pair list below
# preserves the original document order so that the one value that is both
# an old value and a new value elsewhere (77÷6=) is substituted away before
# it is (re)introduced by a later replacement.

$d = $word.ActiveDocument

$replacements = @(
    @("76÷7=", "28÷9="),
    @("79÷4=", "81÷8="),
    @("92÷3=", "31÷5="),
    @("86÷3=", "74÷7="),
    @("39÷7=", "48÷5="),
    @("36÷2=", "96÷3="),
    @("30÷7=", "68÷2="),
    @("13÷5=", "58÷4="),
    @("78÷3=", "41÷4="),
    @("47÷7=", "23÷5="),
    @("20÷9=", "55÷4="),
    @("99÷8=", "50÷3="),
    @("57÷3=", "51÷7="),
    @("62÷9=", "71÷5="),
    @("31÷6=", "43÷4="),
    @("77÷6=", "48÷2="),
    @("69÷8=", "23÷9="),
    @("23÷7=", "46÷4="),
    @("71÷7=", "79÷7="),
    @("55÷8=", "65÷4="),
    @("58÷2=", "23÷8="),
    @("64÷8=", "16÷4="),
    @("18÷6=", "38÷5="),
    @("66÷4=", "77÷6="),
    @("55÷2=", "52÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
